# Auto-generated edit script: update FFXIV Leve profit calculations
# across 7 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW) per upstream data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 52879.15
$ws.Range("I76").Value = 55514.895
$ws.Range("J76").Value = 2800
$ws.Range("K76").Value = 55514.895
$ws.Range("L76").Value = 2800
$ws.Range("M76").Value = -55199.895
$ws.Range("N76").Value = -3430
$ws.Range("H79").Value = 52879.15
$ws.Range("I79").Value = 55514.895
$ws.Range("J79").Value = 2800
$ws.Range("K79").Value = 55514.895
$ws.Range("L79").Value = 2800
$ws.Range("M79").Value = -54422.895
$ws.Range("N79").Value = -4984
$ws.Range("H86").Value = 77001680
$ws.Range("I86").Value = 112899.555
$ws.Range("J86").Value = 250001420
$ws.Range("K86").Value = 112899.555
$ws.Range("L86").Value = 250001420
$ws.Range("M86").Value = -111776.555
$ws.Range("N86").Value = -250003666
$ws.Range("H89").Value = 77001680
$ws.Range("I89").Value = 112899.555
$ws.Range("J89").Value = 250001420
$ws.Range("K89").Value = 564497.7749999999
$ws.Range("L89").Value = 1250007100
$ws.Range("M89").Value = -558881.7749999999
$ws.Range("N89").Value = -1250018332
$ws.Range("H116").Value = 4558.8
$ws.Range("I116").Value = 888
$ws.Range("K116").Value = 888
$ws.Range("M116").Value = 2554
$ws.Range("H132").Value = 1787273.2
$ws.Range("I132").Value = 1906296.8
$ws.Range("J132").Value = 1920
$ws.Range("K132").Value = 5718890.4
$ws.Range("L132").Value = 5760
$ws.Range("M132").Value = -5716360.4
$ws.Range("N132").Value = -10820

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1220.6154
$ws.Range("I2").Value = 627.2857
$ws.Range("J2").Value = 1912.8334
$ws.Range("K2").Value = 627.2857
$ws.Range("L2").Value = 1912.8334
$ws.Range("M2").Value = -514.2857
$ws.Range("N2").Value = -2138.8334
$ws.Range("H5").Value = 74074190
$ws.Range("I5").Value = 23809638
$ws.Range("K5").Value = 23809638
$ws.Range("M5").Value = -23809526
$ws.Range("H116").Value = 1220.6154
$ws.Range("I116").Value = 627.2857
$ws.Range("J116").Value = 1912.8334
$ws.Range("K116").Value = 627.2857
$ws.Range("L116").Value = 1912.8334
$ws.Range("M116").Value = 1666.7143
$ws.Range("N116").Value = -6500.8334

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1220.6154
$ws.Range("I3").Value = 627.2857
$ws.Range("J3").Value = 1912.8334
$ws.Range("K3").Value = 627.2857
$ws.Range("L3").Value = 1912.8334
$ws.Range("M3").Value = -513.2857
$ws.Range("N3").Value = -2140.8334
$ws.Range("H4").Value = 74074190
$ws.Range("I4").Value = 23809638
$ws.Range("K4").Value = 23809638
$ws.Range("M4").Value = -23809523
$ws.Range("H94").Value = 469.7647
$ws.Range("I94").Value = 481.6875
$ws.Range("J94").Value = 279
$ws.Range("K94").Value = 481.6875
$ws.Range("L94").Value = 279
$ws.Range("M94").Value = -30.6875
$ws.Range("N94").Value = -1181

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 165500
$ws.Range("I6").Value = 165500
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 165500
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -165387
$ws.Range("H93").Value = 55555
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 55555
$ws.Range("K93").Value = 0
$ws.Range("L93").ClearContents()
$ws.Range("M93").Value = 55555
$ws.Range("N93").Value = -59299
$ws.Range("H107").Value = 718.5833
$ws.Range("I107").Value = 738.55554
$ws.Range("J107").Value = 688.625
$ws.Range("K107").Value = 738.55554
$ws.Range("L107").Value = 688.625
$ws.Range("M107").Value = 1181.44446
$ws.Range("N107").Value = -4528.625

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 331.65625
$ws.Range("I40").Value = 73.5
$ws.Range("J40").Value = 589.8125
$ws.Range("K40").Value = 294
$ws.Range("L40").Value = 2359.25
$ws.Range("M40").Value = -225
$ws.Range("N40").Value = -2497.25
$ws.Range("H68").Value = 499
$ws.Range("I68").Value = 498.8
$ws.Range("J68").Value = 500
$ws.Range("K68").Value = 1496.4
$ws.Range("L68").Value = 1500
$ws.Range("M68").Value = -685.4000000000001
$ws.Range("N68").Value = -3122
$ws.Range("H71").Value = 499
$ws.Range("I71").Value = 498.8
$ws.Range("J71").Value = 500
$ws.Range("K71").Value = 4489.2
$ws.Range("L71").Value = 4500
$ws.Range("M71").Value = -433.1999999999998
$ws.Range("N71").Value = -12612
$ws.Range("H116").Value = 9095.643
$ws.Range("I116").Value = 12760
$ws.Range("J116").Value = 2499.8
$ws.Range("K116").Value = 38280
$ws.Range("L116").Value = 7499.400000000001
$ws.Range("M116").Value = -34838
$ws.Range("N116").Value = -14383.4
$ws.Range("H121").Value = 1562456.1
$ws.Range("I121").Value = 17092.666
$ws.Range("J121").Value = 2107878.5
$ws.Range("K121").Value = 51277.99800000001
$ws.Range("L121").Value = 6323635.5
$ws.Range("M121").Value = -49967.99800000001
$ws.Range("N121").Value = -6326255.5
$ws.Range("H131").Value = 1607559.1
$ws.Range("J131").Value = 1725284.1
$ws.Range("L131").Value = 5175852.300000001
$ws.Range("N131").Value = -5185932.300000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3093
$ws.Range("J97").Value = 9000
$ws.Range("L97").Value = 9000
$ws.Range("N97").Value = -9992
$ws.Range("H132").Value = 5192.857
$ws.Range("I132").Value = 6841.6665
$ws.Range("K132").Value = 20524.9995
$ws.Range("M132").Value = -17994.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7319.25
$ws.Range("I122").Value = 9910.799999999999
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 29732.4
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -27282.4
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 9317.77
$ws.Range("I132").Value = 11895.223
$ws.Range("J132").Value = 3518.5
$ws.Range("K132").Value = 35685.669
$ws.Range("L132").Value = 10555.5
$ws.Range("M132").Value = -33155.669
$ws.Range("N132").Value = -15615.5
$ws.Range("H136").Value = 5627.0938
$ws.Range("I136").Value = 6465.0835
$ws.Range("J136").Value = 3113.125
$ws.Range("K136").Value = 19395.2505
$ws.Range("L136").Value = 9339.375
$ws.Range("M136").Value = -16845.2505
$ws.Range("N136").Value = -14439.375
